# "Generate Report for Handback"
#
# The localization status report is refreshed after a handback: the
# "Ready for handoff" status becomes "Handed back: in sync with en-US",
# and each language sheet (zh-cn / de-de) gets its "Latest Target File",
# "Latest Handback File" and "Latest Handback DateTime" columns filled
# in (they were blank / placeholder before).

$wb = $excel.ActiveWorkbook

$mdUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/9f4b1357b5ae6711fa491c2bed3d1dcac7447cb1/e2e/1dc4552f-4109-4511-b3c0-35bc5bd32916.md"
$mdName = "1dc4552f-4109-4511-b3c0-35bc5bd32916.md"
$newStatus = "Handed back: in sync with en-US"

# --- Overview sheet: status shown per-language is now "handed back" ---
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = $newStatus
$overview.Range("F2").Value = $newStatus
$overview.Columns.Item(5).ColumnWidth = 29.166666666666668
$overview.Columns.Item(6).ColumnWidth = 29.166666666666668

# --- zh-cn sheet ---
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = $newStatus
$zhcn.Range("J2").Value = "1dc4552f-4109-4511-b3c0-35bc5bd32916.eef01d68377a84301fb70fd3fff53cd96e888a52.zh-cn.xlf"
$zhcn.Range("K2").Value = "2016-11-09 10:53:31"
$zhcn.Hyperlinks.Add($zhcn.Range("I2"), $mdUrl, [System.Type]::Missing, [System.Type]::Missing, $mdName)
$zhcn.Columns.Item(3).ColumnWidth = 29.166666666666668
$zhcn.Columns.Item(9).ColumnWidth = 39.166666666666664
$zhcn.Columns.Item(10).ColumnWidth = 39.166666666666664

# --- de-de sheet ---
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = $newStatus
$dede.Range("J2").Value = "1dc4552f-4109-4511-b3c0-35bc5bd32916.eef01d68377a84301fb70fd3fff53cd96e888a52.de-de.xlf"
$dede.Range("K2").Value = "2016-11-09 10:53:50"
$dede.Hyperlinks.Add($dede.Range("I2"), $mdUrl, [System.Type]::Missing, [System.Type]::Missing, $mdName)
$dede.Columns.Item(3).ColumnWidth = 29.166666666666668
$dede.Columns.Item(9).ColumnWidth = 39.166666666666664
$dede.Columns.Item(10).ColumnWidth = 39.166666666666664
